$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width: target stored width is 15.42578125 characters.
# ColumnWidth assignments are quantized internally to 1/6-character steps,
# so 14.67 is the nearest input that lands on the closest achievable
# stored width (15.5).
$ws.Columns.Item(1).ColumnWidth = 14.67

$values = @(
    0.46475278436945189,
    -0.053651242282285239,
    -0.0089999999506940043,
    -0.011999999985093268,
    -0.0059999999514133151,
    -0.005999999949640511,
    -0.019999999939201984,
    -0.002110630581223738,
    -0.0059999999486377575,
    -0.0059999999482158728,
    -0.004499999949331368,
    -0.0059999999481101796,
    -0.005999999947788659,
    -0.011999999943311579,
    0.022339889787088119,
    -0.0059999999475603971,
    -0.0059999999473037136,
    -0.0089999999450345314,
    -0.0089999999527554664,
    -0.013800769890750075,
    -0.0089999999521808149,
    -0.0089999999521177543,
    -0.081839425583414283,
    -0.041999999924666298,
    -0.041999999924217768,
    -0.0059999999494593226,
    -0.0059999999491204825,
    -0.0059999999479884991,
    -0.011999999942728934,
    -0.019999999936416657,
    -0.014999999939501407,
    -0.020999999934965707,
    -0.0059999999459776632
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
